$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9693716918425304
$ws.Range("J2").Value = 0.9693716918425304
$ws.Range("M2").Value = 2.133443333333334
$ws.Range("N2").Value = 6.40033
$ws.Range("O2").Value = 0.2605947899689859
$ws.Range("P2").Value = 0.2605947899689859
$ws.Range("Q2").Value = 19.53585384330445
$ws.Range("R2").Value = 175.82268458974
$ws.Range("S2").Value = 0.2526132124375848
$ws.Range("T2").Value = 0.2526132124375848

# Row 3
$ws.Range("I3").Value = 0.9693716918425304
$ws.Range("J3").Value = 0.9693716918425304
$ws.Range("O3").Value = 0.5209338844846115
$ws.Range("P3").Value = 0.5209338844846116
$ws.Range("S3").Value = 0.5049785609409492
$ws.Range("T3").Value = 0.5049785609409493

# Row 4
$ws.Range("I4").Value = 0.9693716918425304
$ws.Range("J4").Value = 0.9693716918425304
$ws.Range("M4").Value = 1.788586
$ws.Range("N4").Value = 5.365758
$ws.Range("O4").Value = 0.2184713255464024
$ws.Range("P4").Value = 0.2184713255464024
$ws.Range("Q4").Value = 16.37800926616933
$ws.Range("R4").Value = 147.402083395524
$ws.Range("S4").Value = 0.2117799184639963
$ws.Range("T4").Value = 0.2117799184639964

# Row 5
$ws.Range("G5").Value = 0.2893236666666667
$ws.Range("H5").Value = 0.867971
$ws.Range("I5").Value = 0.03062830815746963
$ws.Range("J5").Value = 0.03062830815746962
$ws.Range("M5").Value = 2.133443333333334
$ws.Range("N5").Value = 6.40033
$ws.Range("O5").Value = 0.2605947899689859
$ws.Range("P5").Value = 0.2605947899689859
$ws.Range("Q5").Value = 0.6172556478255558
$ws.Range("R5").Value = 5.55530083043
$ws.Range("S5").Value = 0.007981577531401177
$ws.Range("T5").Value = 0.007981577531401175

# Row 6
$ws.Range("G6").Value = 0.2893236666666667
$ws.Range("H6").Value = 0.867971
$ws.Range("I6").Value = 0.03062830815746963
$ws.Range("J6").Value = 0.03062830815746962
$ws.Range("O6").Value = 0.5209338844846115
$ws.Range("P6").Value = 0.5209338844846116
$ws.Range("Q6").Value = 1.233905644775556
$ws.Range("R6").Value = 11.10515080298
$ws.Range("S6").Value = 0.01595532354366237
$ws.Range("T6").Value = 0.01595532354366237

# Row 7
$ws.Range("G7").Value = 0.2893236666666667
$ws.Range("H7").Value = 0.867971
$ws.Range("I7").Value = 0.03062830815746963
$ws.Range("J7").Value = 0.03062830815746962
$ws.Range("M7").Value = 1.788586
$ws.Range("N7").Value = 5.365758
$ws.Range("O7").Value = 0.2184713255464024
$ws.Range("P7").Value = 0.2184713255464024
$ws.Range("Q7").Value = 0.5174802596686666
$ws.Range("R7").Value = 4.657322337018
$ws.Range("S7").Value = 0.006691407082406079
$ws.Range("T7").Value = 0.00669140708240608
